# Automatische test-sync: 2025-06-26 21:42:50
# Add a new log entry (row 21) to the "Logs" sheet and bump the
# "Bestelling / Levering" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$subject   = "Bestel je 3 bussen siliconenspray?"
$sender    = "MailMind Test <mailmind.test@zohomail.eu>"
$body      = "He Johan, `nZou je 3 bussen siliconenspray kunnen bestellen voor in de werkplaats?`nThanks,`nMarc`nSent using {0}"
$category  = "Bestelling / Levering"
$reply     = "Bedankt voor je bericht. Ik neem dit z.s.m. in behandeling."
$timestamp = "2025-06-26 21:42:06"

$logs.Range("A21").Value = $subject
$logs.Range("B21").Value = $sender
$logs.Range("C21").Value = $body
$logs.Range("D21").Value = $category
$logs.Range("E21").Value = $reply
$logs.Range("F21").Value = $timestamp
$logs.Range("G21").Value = "Ja"
$logs.Range("H21").Value = "Ja"
$logs.Range("I21").Value = "Nee"

# Writing the multi-line body auto-expands the row height; restore the
# default (no explicit/custom row height), matching every other row.
$logs.Rows.Item(21).AutoFit()

# Extend the conditional-formatting ranges so they cover the new row too.
$logs.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))
$logs.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))
$logs.Range("H2:H20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H21"))
$logs.Range("I2:I20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I21"))

# Bump the dashboard count for "Bestelling / Levering" from 15 to 16.
$dashboard.Range("B2").Value = 16
